$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Update row 2 and row 3 content
$ws.Range("A2").Value = "STUTTGART_MTL"
$ws.Range("B2").Value = "Yes"
$ws.Range("A3").Value = "BOCHUM_MTL_SIM"
$ws.Range("B3").Value = "Yes"

# Remove the now-unused rows 4-8
$ws.Rows("4:8").Delete()

# Match the saved cursor/selection position from the file
$ws.Range("B7").Select()
